$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 7938.727
$ws.Range("I40").Value = 5888.4
$ws.Range("K40").Value = 5888.4
$ws.Range("M40").Value = -5713.4
$ws.Range("H55").Value = 1376.4546
$ws.Range("I55").Value = 833
$ws.Range("J55").Value = 1497.2222
$ws.Range("K55").Value = 833
$ws.Range("L55").Value = 1497.2222
$ws.Range("M55").Value = -619
$ws.Range("N55").Value = -1925.2222
$ws.Range("H95").Value = 21271.143
$ws.Range("J95").Value = 21271.143
$ws.Range("L95").Value = 21271.143
$ws.Range("N95").Value = -26763.143
$ws.Range("H116").Value = 4722.625
$ws.Range("I116").Value = 4726
$ws.Range("J116").Value = 4699
$ws.Range("K116").Value = 4726
$ws.Range("L116").Value = 4699
$ws.Range("M116").Value = -1284
$ws.Range("N116").Value = -11583
$ws.Range("H125").Value = 931
$ws.Range("I125").Value = 931
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 8379
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -5919
$ws.Range("N125").ClearContents()
$ws.Range("H132").Value = 15298.75
$ws.Range("I132").Value = 15569.643
$ws.Range("K132").Value = 46708.929
$ws.Range("M132").Value = -44178.929
$ws.Range("H137").Value = 4583.353
$ws.Range("I137").Value = 881.5
$ws.Range("J137").Value = 9871.714
$ws.Range("K137").Value = 2644.5
$ws.Range("L137").Value = 29615.142
$ws.Range("M137").Value = -94.5
$ws.Range("N137").Value = -34715.142
$ws.Range("H141").Value = 1076.4
$ws.Range("I141").Value = 1076.4
$ws.Range("K141").Value = 3229.2
$ws.Range("M141").Value = 1950.8

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3301.75
$ws.Range("I74").Value = 2738.5557
$ws.Range("J74").Value = 4991.3335
$ws.Range("K74").Value = 2738.5557
$ws.Range("L74").Value = 4991.3335
$ws.Range("M74").Value = -1864.5557
$ws.Range("N74").Value = -6739.3335
$ws.Range("H77").Value = 3301.75
$ws.Range("I77").Value = 2738.5557
$ws.Range("J77").Value = 4991.3335
$ws.Range("K77").Value = 13692.7785
$ws.Range("L77").Value = 24956.6675
$ws.Range("M77").Value = -9324.7785
$ws.Range("N77").Value = -33692.6675

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 27344.555
$ws.Range("I82").Value = 7607.8335
$ws.Range("K82").Value = 7607.8335
$ws.Range("M82").Value = -7224.8335
$ws.Range("H85").Value = 27344.555
$ws.Range("I85").Value = 7607.8335
$ws.Range("K85").Value = 7607.8335
$ws.Range("M85").Value = -6281.8335
$ws.Range("H94").Value = 1799.3334
$ws.Range("I94").Value = 1399.5
$ws.Range("K94").Value = 1399.5
$ws.Range("M94").Value = -948.5
$ws.Range("H95").Value = 20000
$ws.Range("J95").Value = 20000
$ws.Range("L95").Value = 20000
$ws.Range("N95").Value = -25492
$ws.Range("H97").Value = 13249.5
$ws.Range("I97").Value = 13249.5
$ws.Range("K97").Value = 13249.5
$ws.Range("M97").Value = -12258.5
$ws.Range("H99").Value = 974.6667
$ws.Range("I99").Value = 974.6667
$ws.Range("K99").Value = 974.6667
$ws.Range("M99").Value = 523.3333
$ws.Range("H103").Value = 28827
$ws.Range("J103").Value = 28436
$ws.Range("L103").Value = 28436
$ws.Range("N103").Value = -30780
$ws.Range("H106").Value = 58750
$ws.Range("J106").Value = 58750
$ws.Range("L106").Value = 58750
$ws.Range("N106").Value = -61274

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6140.476
$ws.Range("I31").Value = 1615.1111
$ws.Range("J31").Value = 9534.5
$ws.Range("K31").Value = 1615.1111
$ws.Range("L31").Value = 9534.5
$ws.Range("M31").Value = -1320.1111
$ws.Range("N31").Value = -10124.5
$ws.Range("H34").Value = 6140.476
$ws.Range("I34").Value = 1615.1111
$ws.Range("J34").Value = 9534.5
$ws.Range("K34").Value = 1615.1111
$ws.Range("L34").Value = 9534.5
$ws.Range("M34").Value = -1413.1111
$ws.Range("N34").Value = -9938.5
$ws.Range("H58").Value = 3536.9092
$ws.Range("I58").Value = 1490.25
$ws.Range("K58").Value = 1490.25
$ws.Range("M58").Value = -1287.25
$ws.Range("H94").Value = 11004.667
$ws.Range("I94").Value = 4000
$ws.Range("K94").Value = 4000
$ws.Range("M94").Value = -3549
$ws.Range("H99").Value = 3332.5
$ws.Range("I99").Value = 3332.5
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3332.5
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1834.5
$ws.Range("N99").ClearContents()
$ws.Range("H126").Value = 3332.5
$ws.Range("I126").Value = 3332.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 9997.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -7527.5
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 2784.2856
$ws.Range("I132").Value = 2834.6
$ws.Range("J132").Value = 2658.5
$ws.Range("K132").Value = 8503.799999999999
$ws.Range("L132").Value = 7975.5
$ws.Range("M132").Value = -5973.799999999999
$ws.Range("N132").Value = -13035.5
$ws.Range("H136").Value = 3536.9092
$ws.Range("I136").Value = 1490.25
$ws.Range("K136").Value = 4470.75
$ws.Range("M136").Value = -1920.75

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 5283.5
$ws.Range("J36").Value = 7825
$ws.Range("L36").Value = 23475
$ws.Range("N36").Value = -23813
$ws.Range("H38").Value = 445.42105
$ws.Range("I38").Value = 447.6875
$ws.Range("K38").Value = 1343.0625
$ws.Range("M38").Value = -996.0625
$ws.Range("H41").Value = 10000
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 10000
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 30000
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -30676
$ws.Range("H55").Value = 4821.1904
$ws.Range("J55").Value = 5529.8887
$ws.Range("L55").Value = 16589.6661
$ws.Range("N55").Value = -16943.6661
$ws.Range("H59").Value = 665
$ws.Range("I59").Value = 665
$ws.Range("K59").Value = 1995
$ws.Range("M59").Value = -1455
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 39998
$ws.Range("J15").Value = 39998
$ws.Range("L15").Value = 39998
$ws.Range("N15").Value = -40574
$ws.Range("H81").Value = 39998
$ws.Range("J81").Value = 39998
$ws.Range("L81").Value = 39998
$ws.Range("N81").Value = -41994
$ws.Range("H84").Value = 39998
$ws.Range("J84").Value = 39998
$ws.Range("L84").Value = 119994
$ws.Range("N84").Value = -129978
$ws.Range("H97").Value = 832.8570999999999
$ws.Range("I97").Value = 587.5
$ws.Range("J97").Value = 1160
$ws.Range("K97").Value = 587.5
$ws.Range("L97").Value = 1160
$ws.Range("M97").Value = -91.5
$ws.Range("N97").Value = -2152
$ws.Range("H113").Value = 2064.1667
$ws.Range("I113").Value = 1977
$ws.Range("K113").Value = 1977
$ws.Range("M113").Value = 193

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1649.25
$ws.Range("J16").Value = 1599
$ws.Range("L16").Value = 1599
$ws.Range("N16").Value = -1939
$ws.Range("H61").Value = 7698.75
$ws.Range("I61").Value = 6796
$ws.Range("K61").Value = 6796
$ws.Range("M61").Value = -6594
$ws.Range("H93").Value = 1339.8572
$ws.Range("I93").Value = 1462.2222
$ws.Range("K93").Value = 1462.2222
$ws.Range("M93").Value = -214.2221999999999
$ws.Range("H113").Value = 7698.75
$ws.Range("I113").Value = 6796
$ws.Range("K113").Value = 6796
$ws.Range("M113").Value = -4626

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 472.2
$ws.Range("I96").Value = 472.2
$ws.Range("K96").Value = 472.2
$ws.Range("M96").Value = 900.8
$ws.Range("H97").Value = 10000
$ws.Range("J97").Value = 10000
$ws.Range("L97").Value = 10000
$ws.Range("N97").Value = -11982
$ws.Range("H100").Value = 1132.4166
$ws.Range("I100").Value = 1132.4166
$ws.Range("K100").Value = 2264.8332
$ws.Range("M100").Value = -1723.8332
$ws.Range("H122").Value = 1396.0667
$ws.Range("I122").Value = 1233
$ws.Range("K122").Value = 3699
$ws.Range("M122").Value = -1249
